# Update the cryptos price/volume list (GitHub Actions refresh).
# Numeric-looking "Price" strings are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.228.98"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.636.02"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'216.58"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'0.523"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.257"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'20.37"
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").Value = "'0.0850"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.639.97"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "'65.36"
$ws.Range("D16").Value = "27.194.93"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "0.0₃0743"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'218.67"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'6.98"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'2.42"
$ws.Range("E22").Value = "  -6.27%  "
$ws.Range("D23").Value = "'9.08"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "'147.82"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'7.33"
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "'15.70"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "1.340.75"
$ws.Range("E33").Value = "  +5.65%  "
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'0.549"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'0.854"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("E40").Value = "  +1.43%  "
$ws.Range("D41").Value = "'0.804"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").Value = "'64.87"
$ws.Range("E42").Value = "  +4.58%  "
$ws.Range("D43").Value = "1.775.27"
$ws.Range("D44").Value = "'5.27"
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("D45").Value = "'90.92"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'0.814"
$ws.Range("E47").Value = "  +21.93%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.0514"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.0994"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0991"
$ws.Range("E50").Value = "  -7.83%  "
$ws.Range("D51").Value = "'7.58"
$ws.Range("E51").Value = "  -0.68%  "
